$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3732763333333333
$ws.Range("H2").Value = 1.119829
$ws.Range("I2").Value = 0.3554258969843855
$ws.Range("J2").Value = 0.3554258969843855
$ws.Range("M2").Value = 0.110552
$ws.Range("N2").Value = 0.331656
$ws.Range("O2").Value = 0.01126249561724847
$ws.Range("P2").Value = 0.01126249561724847
$ws.Range("Q2").Value = 0.04126644520266667
$ws.Range("R2").Value = 0.371398006824
$ws.Range("S2").Value = 0.004002982607043249
$ws.Range("T2").Value = 0.004002982607043249
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3732763333333333
$ws.Range("H3").Value = 1.119829
$ws.Range("I3").Value = 0.3554258969843855
$ws.Range("J3").Value = 0.3554258969843855
$ws.Range("O3").Value = 0.9181055646724333
$ws.Range("P3").Value = 0.9181055646724334
$ws.Range("Q3").Value = 3.363992694194222
$ws.Range("R3").Value = 30.275934247748
$ws.Range("S3").Value = 0.3263184938500553
$ws.Range("T3").Value = 0.3263184938500554
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3732763333333333
$ws.Range("H4").Value = 1.119829
$ws.Range("I4").Value = 0.3554258969843855
$ws.Range("J4").Value = 0.3554258969843855
$ws.Range("M4").Value = 0.6933189999999999
$ws.Range("N4").Value = 2.079957
$ws.Range("O4").Value = 0.07063193971031816
$ws.Range("P4").Value = 0.07063193971031817
$ws.Range("Q4").Value = 0.2587995741503333
$ws.Range("R4").Value = 2.329196167353
$ws.Range("S4").Value = 0.02510442052728687
$ws.Range("T4").Value = 0.02510442052728687
$ws.Range("G5").Value = 0.668317
$ws.Range("I5").Value = 0.6363574327729865
$ws.Range("J5").Value = 0.6363574327729865
$ws.Range("M5").Value = 0.110552
$ws.Range("N5").Value = 0.331656
$ws.Range("O5").Value = 0.01126249561724847
$ws.Range("P5").Value = 0.01126249561724847
$ws.Range("Q5").Value = 0.073883780984
$ws.Range("R5").Value = 0.6649540288560001
$ws.Range("S5").Value = 0.007166972797609249
$ws.Range("T5").Value = 0.007166972797609251
$ws.Range("G6").Value = 0.668317
$ws.Range("I6").Value = 0.6363574327729865
$ws.Range("J6").Value = 0.6363574327729865
$ws.Range("O6").Value = 0.9181055646724333
$ws.Range("P6").Value = 0.9181055646724334
$ws.Range("Q6").Value = 6.022920031734667
$ws.Range("R6").Value = 54.206280285612
$ws.Range("S6").Value = 0.5842433001495428
$ws.Range("T6").Value = 0.5842433001495428
$ws.Range("G7").Value = 0.668317
$ws.Range("I7").Value = 0.6363574327729865
$ws.Range("J7").Value = 0.6363574327729865
$ws.Range("M7").Value = 0.6933189999999999
$ws.Range("N7").Value = 2.079957
$ws.Range("O7").Value = 0.07063193971031816
$ws.Range("P7").Value = 0.07063193971031817
$ws.Range("Q7").Value = 0.463356874123
$ws.Range("R7").Value = 4.170211867107
$ws.Range("S7").Value = 0.04494715982583442
$ws.Range("T7").Value = 0.04494715982583443
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.008629333333333334
$ws.Range("H8").Value = 0.025888
$ws.Range("I8").Value = 0.008216670242627913
$ws.Range("J8").Value = 0.008216670242627911
$ws.Range("M8").Value = 0.110552
$ws.Range("N8").Value = 0.331656
$ws.Range("O8").Value = 0.01126249561724847
$ws.Range("P8").Value = 0.01126249561724847
$ws.Range("Q8").Value = 0.0009539900586666667
$ws.Range("R8").Value = 0.008585910528
$ws.Range("S8").Value = 0.00009254021259597281
$ws.Range("T8").Value = 0.0000925402125959728
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.008629333333333334
$ws.Range("H9").Value = 0.025888
$ws.Range("I9").Value = 0.008216670242627913
$ws.Range("J9").Value = 0.008216670242627911
$ws.Range("O9").Value = 0.9181055646724333
$ws.Range("P9").Value = 0.9181055646724334
$ws.Range("Q9").Value = 0.07776816180622223
$ws.Range("R9").Value = 0.699913456256
$ws.Range("S9").Value = 0.00754377067283508
$ws.Range("T9").Value = 0.007543770672835079
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.008629333333333334
$ws.Range("H10").Value = 0.025888
$ws.Range("I10").Value = 0.008216670242627913
$ws.Range("J10").Value = 0.008216670242627911
$ws.Range("M10").Value = 0.6933189999999999
$ws.Range("N10").Value = 2.079957
$ws.Range("O10").Value = 0.07063193971031816
$ws.Range("P10").Value = 0.07063193971031817
$ws.Range("Q10").Value = 0.005982880757333334
$ws.Range("R10").Value = 0.053845926816
$ws.Range("S10").Value = 0.00058035935719686
$ws.Range("T10").Value = 0.00058035935719686
